$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 204.5
$ws.Range("I12").Value = 209
$ws.Range("K12").Value = 209
$ws.Range("M12").Value = -39
$ws.Range("H38").Value = 1263.35
$ws.Range("I38").Value = 114.36364
$ws.Range("J38").Value = 2667.6667
$ws.Range("K38").Value = 343.09092
$ws.Range("L38").Value = 8003.000100000001
$ws.Range("M38").Value = 28.90908000000002
$ws.Range("N38").Value = -8747.000100000001
$ws.Range("H58").Value = 605
$ws.Range("I58").Value = 47
$ws.Range("K58").Value = 141
$ws.Range("M58").Value = 9
$ws.Range("H87").Value = 40260.332
$ws.Range("J87").Value = 40260.332
$ws.Range("L87").Value = 40260.332
$ws.Range("N87").Value = -42756.332
$ws.Range("H90").Value = 40260.332
$ws.Range("J90").Value = 40260.332
$ws.Range("L90").Value = 120780.996
$ws.Range("N90").Value = -133260.996
$ws.Range("H132").Value = 2954.4363
$ws.Range("I132").Value = 2104.1777
$ws.Range("J132").Value = 6780.6
$ws.Range("K132").Value = 6312.533100000001
$ws.Range("L132").Value = 20341.8
$ws.Range("M132").Value = -3782.533100000001
$ws.Range("N132").Value = -25401.8
$ws.Range("H133").Value = 39997.777
$ws.Range("J133").Value = 39997.777
$ws.Range("L133").Value = 39997.777
$ws.Range("N133").Value = -50117.777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1022.2
$ws.Range("I110").Value = 1027.75
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 1027.75
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 1017.25
$ws.Range("N110").Value = -5090
$ws.Range("H132").Value = 2269.6155
$ws.Range("I132").Value = 866.9231
$ws.Range("K132").Value = 2600.7693
$ws.Range("M132").Value = -70.76929999999993

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1728.4419
$ws.Range("I134").Value = 1387.6562
$ws.Range("J134").Value = 2719.818
$ws.Range("K134").Value = 4162.9686
$ws.Range("L134").Value = 8159.454000000001
$ws.Range("M134").Value = -1627.9686
$ws.Range("N134").Value = -13229.454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2810
$ws.Range("I62").Value = 2475
$ws.Range("J62").Value = 4150
$ws.Range("K62").Value = 2475
$ws.Range("L62").Value = 4150
$ws.Range("M62").Value = -1851
$ws.Range("N62").Value = -5398
$ws.Range("H65").Value = 2810
$ws.Range("I65").Value = 2475
$ws.Range("J65").Value = 4150
$ws.Range("K65").Value = 12375
$ws.Range("L65").Value = 20750
$ws.Range("M65").Value = -9255
$ws.Range("N65").Value = -26990
$ws.Range("H86").Value = 4504.56
$ws.Range("I86").Value = 4884.9473
$ws.Range("J86").Value = 3300
$ws.Range("K86").Value = 4884.9473
$ws.Range("L86").Value = 3300
$ws.Range("M86").Value = -3761.9473
$ws.Range("N86").Value = -5546
$ws.Range("H89").Value = 4504.56
$ws.Range("I89").Value = 4884.9473
$ws.Range("J89").Value = 3300
$ws.Range("K89").Value = 24424.7365
$ws.Range("L89").Value = 16500
$ws.Range("M89").Value = -18808.7365
$ws.Range("N89").Value = -27732

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 1920
$ws.Range("I124").Value = 1560
$ws.Range("K124").Value = 4680
$ws.Range("M124").Value = 230
$ws.Range("H140").Value = 2864.9722
$ws.Range("I140").Value = 1894.037
$ws.Range("J140").Value = 5777.778
$ws.Range("K140").Value = 5682.111
$ws.Range("L140").Value = 17333.334
$ws.Range("M140").Value = -502.1109999999999
$ws.Range("N140").Value = -27693.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 40.526318
$ws.Range("I2").Value = 54.25
$ws.Range("J2").Value = 17
$ws.Range("K2").Value = 54.25
$ws.Range("L2").Value = 17
$ws.Range("M2").Value = 58.75
$ws.Range("N2").Value = -243
$ws.Range("H43").Value = 9958.5
$ws.Range("I43").Value = 9958.5
$ws.Range("K43").Value = 9958.5
$ws.Range("M43").Value = -9807.5
$ws.Range("H46").Value = 1541
$ws.Range("I46").Value = 1541
$ws.Range("K46").Value = 1541
$ws.Range("M46").Value = -1385
$ws.Range("H57").Value = 8816.143
$ws.Range("I57").Value = 2000
$ws.Range("J57").Value = 9952.166999999999
$ws.Range("K57").Value = 2000
$ws.Range("L57").Value = 9952.166999999999
$ws.Range("M57").Value = -1180
$ws.Range("N57").Value = -11592.167
$ws.Range("H80").Value = 2889.5908
$ws.Range("I80").Value = 2718.3333
$ws.Range("J80").Value = 3008.1538
$ws.Range("K80").Value = 2718.3333
$ws.Range("L80").Value = 3008.1538
$ws.Range("M80").Value = -1720.3333
$ws.Range("N80").Value = -5004.1538
$ws.Range("H83").Value = 2889.5908
$ws.Range("I83").Value = 2718.3333
$ws.Range("J83").Value = 3008.1538
$ws.Range("K83").Value = 13591.6665
$ws.Range("L83").Value = 15040.769
$ws.Range("M83").Value = -8599.666499999999
$ws.Range("N83").Value = -25024.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1345.3636
$ws.Range("I61").Value = 1185.5714
$ws.Range("J61").Value = 1625
$ws.Range("K61").Value = 1185.5714
$ws.Range("L61").Value = 1625
$ws.Range("M61").Value = -983.5714
$ws.Range("N61").Value = -2029
$ws.Range("H82").Value = 2335.1052
$ws.Range("I82").Value = 2798
$ws.Range("J82").Value = 1918.5
$ws.Range("K82").Value = 2798
$ws.Range("L82").Value = 1918.5
$ws.Range("M82").Value = -2437
$ws.Range("N82").Value = -2640.5
$ws.Range("H85").Value = 2335.1052
$ws.Range("I85").Value = 2798
$ws.Range("J85").Value = 1918.5
$ws.Range("K85").Value = 2798
$ws.Range("L85").Value = 1918.5
$ws.Range("M85").Value = -1550
$ws.Range("N85").Value = -4414.5
$ws.Range("H113").Value = 1345.3636
$ws.Range("I113").Value = 1185.5714
$ws.Range("J113").Value = 1625
$ws.Range("K113").Value = 1185.5714
$ws.Range("L113").Value = 1625
$ws.Range("M113").Value = 984.4286
$ws.Range("N113").Value = -5965

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 874.93335
$ws.Range("I113").Value = 810.6
$ws.Range("J113").Value = 1003.6
$ws.Range("K113").Value = 2431.8
$ws.Range("L113").Value = 3010.8
$ws.Range("M113").Value = -261.8000000000002
$ws.Range("N113").Value = -7350.8
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
